$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.271.80'
$ws.Range("E2").Value = '  -1.05%  '
$ws.Range("D3").Value = '2.517.19'
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.56'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.36%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.515'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.24%  '
$ws.Range("D9").Value = '2.515.02'
$ws.Range("E9").Value = '  -0.14%  '
$ws.Range("E10").Value = '  +1.64%  '
$ws.Range("E11").Value = '  -0.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.356'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.08%  '
$ws.Range("E13").Value = '  +1.95%  '
$ws.Range("D14").Value = '2.980.52'
$ws.Range("E14").Value = '  +0.03%  '
$ws.Range("D15").Value = '69.098.30'
$ws.Range("E15").Value = '  -1.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000175'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.49%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.79'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("D18").Value = '2.515.38'
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.32'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.78%  '
$ws.Range("E20").Value = '  +1.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '347.60'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.46%  '
$ws.Range("E22").Value = '  -0.43%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.05%  '
$ws.Range("E24").Value = '  +0.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.46'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.40%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.95'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.89%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.86'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.40%  '
$ws.Range("D28").Value = '2.649.16'
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("D30").Value = '0.0₃0889'
$ws.Range("E30").Value = '  -1.90%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.80'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '460.59'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.72%  '
$ws.Range("E33").Value = '  -3.71%  '
$ws.Range("E34").Value = '  -1.28%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("E36").Value = '  +0.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '157.39'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.52'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.12%  '
$ws.Range("E40").Value = '  -0.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.318'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.38%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.70'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.60'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '38.13'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.12'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.78%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.21'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.62%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '141.82'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("E48").Value = '  -0.35%  '
$ws.Range("E49").Value = '  -1.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0731'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.579'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.09%  '
